$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 239, pushing existing rows 239:245 down to 240:246
$ws.Rows.Item(239).Insert()

# Populate the new row 239 with the new weekly data point
$ws.Cells.Item(239, 1).Value = 10
$ws.Cells.Item(239, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(239, 3).Value = "La Araucanía"
$ws.Cells.Item(239, 4).Value = 44509
$ws.Cells.Item(239, 5).Value = 9
$ws.Cells.Item(239, 6).Value = 100112037
$ws.Cells.Item(239, 7).Value = "Cebollín"
$ws.Cells.Item(239, 8).Value = "Sin especificar"
$ws.Cells.Item(239, 9).Value = "Primera"
$ws.Cells.Item(239, 10).Value = 30
$ws.Cells.Item(239, 11).Value = 8000
$ws.Cells.Item(239, 12).Value = 8000
$ws.Cells.Item(239, 13).Value = 8000
$ws.Cells.Item(239, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(239, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(239, 16).Value = 667
$ws.Cells.Item(239, 17).Value = 12
$ws.Cells.Item(239, 18).Value = "Hortaliza"
